# adds RVCAT to FN2 species code conversion
# Append a new lookup row (Bloat / 204 / "094") to the SPC_SHORT -> SPECIES/SPC
# conversion table on Sheet1, then leave the view scrolled/selected the way
# the author left it (viewport near row 12, active cell C41 just below the
# newly-entered data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 40: A = species short code, B = RVCAT numeric code,
# C = FN2 species code stored as text (leading-zero code "094", so it needs
# the leading apostrophe to keep it text with a quote-prefix, matching the
# existing "091"/"093"/"081" ... entries already in column C).
$ws.Range("A40").Value = "Bloat"
$ws.Range("B40").Value = 204
$ws.Range("C40").Value = "'094"

# Match the author's final scroll position / selection in the saved file.
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C41").Select()
